$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L (12) to make room for spectraltype_esphs_dr2.
# This shifts the old L:Y (dr3 id/astrometry columns) one to the right, to M:Z.
$ws.Columns.Item(12).Insert()

# New header for the inserted column, and for the appended column at the end (AA).
$ws.Range("L1").Value = "spectraltype_esphs_dr2"
$ws.Range("AA1").Value = "spectraltype_esphs_dr3"

# Column widths: the new column (L) is wider than its neighbours, and so is the
# brand-new trailing column (AA). COM's ColumnWidth stores 5/6 narrower than the
# raw OOXML "width" attribute ends up holding, so compensate by that offset.
$widthOffset = 5 / 6
$ws.Columns.Item(12).ColumnWidth = 24 - $widthOffset
$ws.Columns.Item(27).ColumnWidth = 24 - $widthOffset

# Spectral-type values (Gaia ESP-HS spectral class) for dr2 (col L) and dr3 (col AA).
# Rows that never had a full dr3 cross-match (only L/M populated before the insert)
# only get the dr2 value; rows with no dr2 match at all (the "duplicate slot" rows
# 6-9 and 18-21) get neither.
$l2 = @{
  2="K"; 3="K"; 4="K"; 5="K";
  10="K"; 11="K"; 12="K"; 13="K";
  14="G"; 15="G"; 16="K"; 17="K";
  22="K"; 23="K"; 24="K"; 25="K";
  26="K"; 27="K"
}
$aa3 = @{
  2="K"; 5="K"; 11="K"; 13="K";
  15="G"; 16="K"; 22="K"; 24="K"; 27="K"
}

foreach ($row in $l2.Keys) {
    $ws.Cells.Item($row, 12).Value = $l2[$row]
}
foreach ($row in $aa3.Keys) {
    $ws.Cells.Item($row, 27).Value = $aa3[$row]
}

Write-Output "done"
